# Update the "Avverkningsanmälningar" sheet:
#  - Remove the two newest entries that used to sit at the top (rows 2-3:
#    "A 57664-2023" and "A 57619-2023"), which shifts every remaining
#    record up by two rows.
#  - Remove the three newest entries that used to sit at the bottom
#    (originally rows 29-31: "A 57804-2023", "A 57805-2023", "A 57807-2023").
#  - Bump the "Förändrad" (last-changed) date for every remaining record
#    from 2023-11-20 (45250) to 2023-11-21 (45251).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the newest two rows at the top of the data (original rows 2 and 3).
$ws.Rows("2:3").Delete()

# Delete the newest three rows that were originally at the bottom
# (originally rows 29-31); after the deletion above they now sit at 27-29.
$ws.Rows("27:29").Delete()

# All remaining data rows are now rows 2 through 26. Update the
# "Förändrad" column (C) for each of them to the new date serial 45251
# (2023-11-21).
$ws.Range("C2:C26").Value = 45251
